{"js": "// Office.js (Word JavaScript API) script\n// Applies a series of small wording/typo fixes throughout the document by\n// locating each old phrase with Range.search() and replacing it in place\n// with Range.insertText(..., \"Replace\"). Using generous, uniquely-matching\n// context for every search keeps each edit targeted at the correct spot\n// even though several of the words involved (e.g. \"nh\u00f3m\") repeat elsewhere\n// in the document.\n\nconst edits = [\n  [\n    \"Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , M\u1ed9t v\u1ea5n \u0111\u1ec1 r\u1ea5t ph\u1ed5 bi\u1ebfn v\u00e0 th\u01b0\u1eddng xuy\u00ean g\u1eb7p ph\u1ea3i \u0111\u1ed1i v\u1edbi n\",\n    \"Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , m\u1ed9t v\u1ea5n \u0111\u1ec1 r\u1ea5t ph\u1ed5 bi\u1ebfn v\u00e0 th\u01b0\u1eddng xuy\u00ean g\u1eb7p ph\u1ea3i \u0111\u1ed1i v\u1edbi n\",\n  ],\n  [\n    \"Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n th\u00e0nh c\u00f4ng vi\u1ec7c , th\u00ec h\u1ecd th\u00eam ph\u1ea7n vi\u1ec7c\",\n    \"Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n th\u00e0nh c\u00f4ng vi\u1ec7c th\u00ec h\u1ecd th\u00eam ph\u1ea7n vi\u1ec7c\",\n  ],\n  [\n    \"Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1eddng ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m hi\u1ec7n nay , v\u00e0 n\u00f3 tr\u1edf th\u00e0nh m\u1ed9t trong nh\u1eefng v\u1ea5n \u0111\u1ec1 r\u1ea5t nh\u1ee9c nh\u00f3i v\u00e0 kh\u00f3 ch\u1ecbu trong giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\",\n    \"Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1eddng ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m hi\u1ec7n nay v\u00e0 n\u00f3 tr\u1edf th\u00e0nh m\u1ed9t trong nh\u1eefng v\u1ea5n \u0111\u1ec1 r\u1ea5t nh\u1ee9c nh\u00f3i v\u00e0 kh\u00f3 ch\u1ecbu trong giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\",\n  ],\n  [\n    \"\u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00f3m ph\u00e1t hi\u1ec7n s\u1edbm c\u00e1c v\u1ea5n \u0111\u1ec1\",\n    \"\u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n ph\u00e1t hi\u1ec7n s\u1edbm c\u00e1c v\u1ea5n \u0111\u1ec1\",\n  ],\n  [\n    \"cung c\u1ea5p cho c\u00e1c nh\u00f3m li\u00ean quan v\u1ec1 k\u1ebft qu\u1ea3 x\u00e2y d\u1ef1ng v\u00e0 th\u1eed nghi\u1ec7m\",\n    \"cung c\u1ea5p cho c\u00e1c b\u00ean li\u00ean quan v\u1ec1 k\u1ebft qu\u1ea3 x\u00e2y d\u1ef1ng v\u00e0 th\u1eed nghi\u1ec7m\",\n  ],\n  [\n    \"\u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda b\u1ea5t k\u1ef3 l\u1ed7i n\u00e0o h\u1ecd \",\n    \"\u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda k\u1ecbp th\u1eddi c\u00e1c l\u1ed7i h\u1ecd \",\n  ],\n  [\n    \"th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb ( repository ) \",\n    \"th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb ( shared repository ) .\",\n  ],\n  [\n    \"Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y d\u1ef1ng l\u1ea1i th\u00f4ng qua build server .\",\n    \"Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y d\u1ef1ng l\u1ea1i th\u00f4ng qua m\u00e1y ch\u1ee7 .\",\n  ],\n  [\n    \"Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch hang khi kh\u00f4ng c\u00f2n l\u1ed7i n\u1eefa .\",\n    \"Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch h\u00e0ng khi kh\u00f4ng c\u00f2n l\u1ed7i n\u1eefa .\",\n  ],\n  [\n    \"Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v\u00e0 c\u1ee5 th\u1ec3 cho t\u1eebng giai \u0111o\u1ea1n .\",\n    \"Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v\u00e0 c\u1ee5 th\u1ec3 cho t\u1eebng giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n .\",\n  ],\n  [\n    \"N\u00e2ng c\u00e1o k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\",\n    \"N\u00e2ng cao k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\",\n  ],\n  [\n    \"c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho Ci ,..\",\n    \"c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho CI ,..\",\n  ],\n  [\n    \"gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh T\u00edch h\u1ee3p li\u00ean t\u1ee5c theo \",\n    \"gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh t\u00edch h\u1ee3p li\u00ean t\u1ee5c theo \",\n  ],\n  [\n    \"b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k\u1ef9 s\u01b0 ph\u1ea7n m\u1ec1m v\u1edbi \u00fd \u0111\u1ecbnh h\u01b0\u1edbng t\u1edbi s\u1ef1 \u0111a d\u1ea1ng h\u00f3a cho Jenkins . Hi\u1ec7n t\u1ea1i v\u1edbi h\u01a1n 1000+ plugins , Jenkins c\u00f3 th\u1ec3 t\u00edch h\u1ee3p v\u1edbi g\u1ea7n h\u1ebft c\u00e1c c\u00f4ng c\u1ee5 v\u00e0 n\u1ec1n t\u1ea3ng hi\u1ec7n nay .\",\n    \"b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k\u1ef9 s\u01b0 ph\u1ea7n m\u1ec1m v\u1edbi \u00fd \u0111\u1ecbnh h\u01b0\u1edbng t\u1edbi s\u1ef1 \u0111a d\u1ea1ng h\u00f3a cho Jenkins . Hi\u1ec7n t\u1ea1i v\u1edbi h\u01a1n 1000+ plugins v\u00e0 \u0111ang t\u0103ng  , Jenkins tr\u1edf th\u00e0nh m\u1ed9t c\u00f4ng c\u1ee5 \u0111\u00e1ng ch\u00fa \u00fd cho c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m \u0111\u1ec3 \u00e1p d\u1ee5ng th\u1ef1c ti\u1ec5n CI c\u1ee7a h\u1ecd .\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of edits) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Always take the first (and expected only) match \u2014 every search string\n  // above was chosen to be unique within the document.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# PowerShell / Word COM interop script\n# Applies a series of small wording/typo fixes throughout the document using\n# Find.Execute(..., Replace:=wdReplaceOne) against uniquely-matching, generous\n# search phrases so each fix lands on the correct occurrence even where a word\n# (e.g. \"nhom\") repeats elsewhere in the document.\n\n$d = $word.ActiveDocument\n\n# Edit 1: \"Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , M\u1ed9t...\" -> \"Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , m\u1ed9t...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , M\u1ed9t v\u1ea5n \u0111\u1ec1 r\u1ea5t ph\u1ed5 bi\u1ebfn v\u00e0 th\u01b0\u1eddng xuy\u00ean g\u1eb7p ph\u1ea3i \u0111\u1ed1i v\u1edbi n\", $true, $false, $false, $false, $false, $true, 1, $false, \"Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , m\u1ed9t v\u1ea5n \u0111\u1ec1 r\u1ea5t ph\u1ed5 bi\u1ebfn v\u00e0 th\u01b0\u1eddng xuy\u00ean g\u1eb7p ph\u1ea3i \u0111\u1ed1i v\u1edbi n\", 1)\nif (-not $found) { throw \"Edit 1: text not found: Trong m\u00f4i tr\u01b0\u1eddng t\u00edch h\u1ee3p li\u00ean t\u1ee5c , M\u1ed9t v\u1ea5n \u0111\u1ec1 r\u1ea5t ph\u1ed5 bi\u1ebfn v\u00e0 th\u01b0\u1eddng xuy\u00ean g\u1eb7p ph\u1ea3i \u0111\u1ed1i v\u1edbi n\" }\n\n# Edit 2: \"Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n ...\" -> \"Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n ...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n th\u00e0nh c\u00f4ng vi\u1ec7c , th\u00ec h\u1ecd th\u00eam ph\u1ea7n vi\u1ec7c\", $true, $false, $false, $false, $false, $true, 1, $false, \"Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n th\u00e0nh c\u00f4ng vi\u1ec7c th\u00ec h\u1ecd th\u00eam ph\u1ea7n vi\u1ec7c\", 1)\nif (-not $found) { throw \"Edit 2: text not found: Sau khi m\u1ed9t nh\u00e0 ph\u00e1t tri\u1ec3n , k\u1ef9 s\u01b0 ho\u00e0n th\u00e0nh c\u00f4ng vi\u1ec7c , th\u00ec h\u1ecd th\u00eam ph\u1ea7n vi\u1ec7c\" }\n\n# Edit 3: \"Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1edd...\" -> \"Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1edd...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1eddng ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m hi\u1ec7n nay , v\u00e0 n\u00f3 tr\u1edf th\u00e0nh m\u1ed9t trong nh\u1eefng v\u1ea5n \u0111\u1ec1 r\u1ea5t nh\u1ee9c nh\u00f3i v\u00e0 kh\u00f3 ch\u1ecbu trong giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\", $true, $false, $false, $false, $false, $true, 1, $false, \"Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1eddng ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m hi\u1ec7n nay v\u00e0 n\u00f3 tr\u1edf th\u00e0nh m\u1ed9t trong nh\u1eefng v\u1ea5n \u0111\u1ec1 r\u1ea5t nh\u1ee9c nh\u00f3i v\u00e0 kh\u00f3 ch\u1ecbu trong giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\", 1)\nif (-not $found) { throw \"Edit 3: text not found: Vi\u1ec7c n\u00e0y x\u1ea3y ra r\u1ea5t nhi\u1ec1u trong m\u00f4i tr\u01b0\u1eddng ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m hi\u1ec7n nay , v\u00e0 n\u00f3 tr\u1edf th\u00e0nh m\u1ed9t trong nh\u1eefng v\u1ea5n \u0111\u1ec1 r\u1ea5t nh\u1ee9c nh\u00f3i v\u00e0 kh\u00f3 ch\u1ecbu trong giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\" }\n\n# Edit 4: \"\u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00f3m ph\u00e1t hi\u1ec7n s\u1edbm...\" -> \"\u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n ph\u00e1...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00f3m ph\u00e1t hi\u1ec7n s\u1edbm c\u00e1c v\u1ea5n \u0111\u1ec1\", $true, $false, $false, $false, $false, $true, 1, $false, \"\u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n ph\u00e1t hi\u1ec7n s\u1edbm c\u00e1c v\u1ea5n \u0111\u1ec1\", 1)\nif (-not $found) { throw \"Edit 4: text not found: \u0110i\u1ec1u n\u00e0y cho ph\u00e9p c\u00e1c nh\u00f3m ph\u00e1t hi\u1ec7n s\u1edbm c\u00e1c v\u1ea5n \u0111\u1ec1\" }\n\n# Edit 5: \"cung c\u1ea5p cho c\u00e1c nh\u00f3m li\u00ean quan v\u1ec1 k\u1ebft q...\" -> \"cung c\u1ea5p cho c\u00e1c b\u00ean li\u00ean quan v\u1ec1 k\u1ebft qu...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"cung c\u1ea5p cho c\u00e1c nh\u00f3m li\u00ean quan v\u1ec1 k\u1ebft qu\u1ea3 x\u00e2y d\u1ef1ng v\u00e0 th\u1eed nghi\u1ec7m\", $true, $false, $false, $false, $false, $true, 1, $false, \"cung c\u1ea5p cho c\u00e1c b\u00ean li\u00ean quan v\u1ec1 k\u1ebft qu\u1ea3 x\u00e2y d\u1ef1ng v\u00e0 th\u1eed nghi\u1ec7m\", 1)\nif (-not $found) { throw \"Edit 5: text not found: cung c\u1ea5p cho c\u00e1c nh\u00f3m li\u00ean quan v\u1ec1 k\u1ebft qu\u1ea3 x\u00e2y d\u1ef1ng v\u00e0 th\u1eed nghi\u1ec7m\" }\n\n# Edit 6: \"\u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda b\u1ea5t k\u1ef3 l\u1ed7i n\u00e0o h\u1ecd ...\" -> \"\u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda k\u1ecbp th\u1eddi c\u00e1c l\u1ed7i h\u1ecd ...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"\u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda b\u1ea5t k\u1ef3 l\u1ed7i n\u00e0o h\u1ecd \", $true, $false, $false, $false, $false, $true, 1, $false, \"\u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda k\u1ecbp th\u1eddi c\u00e1c l\u1ed7i h\u1ecd \", 1)\nif (-not $found) { throw \"Edit 6: text not found: \u0111\u1ec3 h\u1ecd c\u00f3 th\u1ec3 s\u1eeda b\u1ea5t k\u1ef3 l\u1ed7i n\u00e0o h\u1ecd \" }\n\n# Edit 7: \"th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb (...\" -> \"th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb (...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb ( repository ) \", $true, $false, $false, $false, $false, $true, 1, $false, \"th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb ( shared repository ) .\", 1)\nif (-not $found) { throw \"Edit 7: text not found: th\u01b0\u1eddng l\u00e0 m\u1ed9t kho l\u01b0u tr\u1eef \u0111\u01b0\u1ee3c chia s\u1ebb ( repository ) \" }\n\n# Edit 8: \"Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y ...\" -> \"Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y ...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y d\u1ef1ng l\u1ea1i th\u00f4ng qua build server .\", $true, $false, $false, $false, $false, $true, 1, $false, \"Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y d\u1ef1ng l\u1ea1i th\u00f4ng qua m\u00e1y ch\u1ee7 .\", 1)\nif (-not $found) { throw \"Edit 8: text not found: Khi code c\u00f3 s\u1ef1 thay \u0111\u1ed5i th\u00ec s\u1ebd \u0111\u01b0\u1ee3c x\u00e2y d\u1ef1ng l\u1ea1i th\u00f4ng qua build server .\" }\n\n# Edit 9: \"Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch...\" -> \"Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch hang khi kh\u00f4ng c\u00f2n l\u1ed7i n\u1eefa .\", $true, $false, $false, $false, $false, $true, 1, $false, \"Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch h\u00e0ng khi kh\u00f4ng c\u00f2n l\u1ed7i n\u1eefa .\", 1)\nif (-not $found) { throw \"Edit 9: text not found: Ph\u00e1t h\u00e0nh phi\u00ean b\u1ea3n ho\u00e0n ch\u1ec9nh cho kh\u00e1ch hang khi kh\u00f4ng c\u00f2n l\u1ed7i n\u1eefa .\" }\n\n# Edit 10: \"Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v...\" -> \"Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v\u00e0 c\u1ee5 th\u1ec3 cho t\u1eebng giai \u0111o\u1ea1n .\", $true, $false, $false, $false, $false, $true, 1, $false, \"Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v\u00e0 c\u1ee5 th\u1ec3 cho t\u1eebng giai \u0111o\u1ea1n ph\u00e1t tri\u1ec3n .\", 1)\nif (-not $found) { throw \"Edit 10: text not found: Cung c\u1ea5p c\u00e1i nh\u00ecn xuy\u00ean su\u1ed1t t\u1ed5ng quan v\u00e0 c\u1ee5 th\u1ec3 cho t\u1eebng giai \u0111o\u1ea1n .\" }\n\n# Edit 11: \"N\u00e2ng c\u00e1o k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean p...\" -> \"N\u00e2ng cao k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean p...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"N\u00e2ng c\u00e1o k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\", $true, $false, $false, $false, $false, $true, 1, $false, \"N\u00e2ng cao k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\", 1)\nif (-not $found) { throw \"Edit 11: text not found: N\u00e2ng c\u00e1o k\u1ef9 n\u0103ng c\u1ee7a \u0111\u1ed9i ng\u0169 nh\u00e2n vi\u00ean ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m .\" }\n\n# Edit 12: \"c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho Ci ,.....\" -> \"c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho CI ,.....\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho Ci ,..\", $true, $false, $false, $false, $false, $true, 1, $false, \"c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho CI ,..\", 1)\nif (-not $found) { throw \"Edit 12: text not found: c\u00e1c c\u00f4ng c\u1ee5 h\u1ed7 tr\u1ee3 cho Ci ,..\" }\n\n# Edit 13: \"gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh T\u00edch h\u1ee3p li\u00ean t\u1ee5...\" -> \"gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh t\u00edch h\u1ee3p li\u00ean t\u1ee5...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh T\u00edch h\u1ee3p li\u00ean t\u1ee5c theo \", $true, $false, $false, $false, $false, $true, 1, $false, \"gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh t\u00edch h\u1ee3p li\u00ean t\u1ee5c theo \", 1)\nif (-not $found) { throw \"Edit 13: text not found: gi\u00fap \u0111\u1ea1t \u0111\u01b0\u1ee3c quy tr\u00ecnh T\u00edch h\u1ee3p li\u00ean t\u1ee5c theo \" }\n\n# Edit 14: \"b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k...\" -> \"b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k...\"\n$rng = $d.Content\n$found = $rng.Find.Execute(\"b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k\u1ef9 s\u01b0 ph\u1ea7n m\u1ec1m v\u1edbi \u00fd \u0111\u1ecbnh h\u01b0\u1edbng t\u1edbi s\u1ef1 \u0111a d\u1ea1ng h\u00f3a cho Jenkins . Hi\u1ec7n t\u1ea1i v\u1edbi h\u01a1n 1000+ plugins , Jenkins c\u00f3 th\u1ec3 t\u00edch h\u1ee3p v\u1edbi g\u1ea7n h\u1ebft c\u00e1c c\u00f4ng c\u1ee5 v\u00e0 n\u1ec1n t\u1ea3ng hi\u1ec7n nay .\", $true, $false, $false, $false, $false, $true, 1, $false, \"b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k\u1ef9 s\u01b0 ph\u1ea7n m\u1ec1m v\u1edbi \u00fd \u0111\u1ecbnh h\u01b0\u1edbng t\u1edbi s\u1ef1 \u0111a d\u1ea1ng h\u00f3a cho Jenkins . Hi\u1ec7n t\u1ea1i v\u1edbi h\u01a1n 1000+ plugins v\u00e0 \u0111ang t\u0103ng  , Jenkins tr\u1edf th\u00e0nh m\u1ed9t c\u00f4ng c\u1ee5 \u0111\u00e1ng ch\u00fa \u00fd cho c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n ph\u1ea7n m\u1ec1m \u0111\u1ec3 \u00e1p d\u1ee5ng th\u1ef1c ti\u1ec5n CI c\u1ee7a h\u1ecd .\", 1)\nif (-not $found) { throw \"Edit 14: text not found: b\u1edfi c\u1ed9ng \u0111\u1ed3ng c\u00e1c nh\u00e0 ph\u00e1t tri\u1ec3n , c\u00e1c k\u1ef9 s\u01b0 ph\u1ea7n m\u1ec1m v\u1edbi \u00fd \u0111\u1ecbnh h\u01b0\u1edbng t\u1edbi s\u1ef1 \u0111a d\u1ea1ng h\u00f3a cho Jenkins . Hi\u1ec7n t\u1ea1i v\u1edbi h\u01a1n 1000+ plugins , Jenkins c\u00f3 th\u1ec3 t\u00edch h\u1ee3p v\u1edbi g\u1ea7n h\u1ebft c\u00e1c c\u00f4ng c\u1ee5 v\u00e0 n\u1ec1n t\u1ea3ng hi\u1ec7n nay .\" }\n"}
